# payment functionality completed by Anil
# - adds BCD00xx contract-id values (replacing the old "=Cn" formula on Sheet1!A2:A4)
# - adds a new "M" column (amount = 1000) on Sheet1 rows 2-4
# - adds a new wide column (R) on Sheet1
# - updates the active selection / scroll position on Sheet1
# - Sheet2's A2/A3 formulas (=Sheet1!A2) recalc automatically to the new cached value

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Sheet1: replace the "=Cn" formulas in column A with literal ContractId values ---
# Order matters: new shared strings are appended in first-use order, and the
# target file appends them as BCD0019, BCD0013, BCD0015 (in that order), while
# they land on rows A4, A2, A3 respectively.
$ws1.Range("A4").Value = "BCD0019"
$ws1.Range("A2").Value = "BCD0013"
$ws1.Range("A3").Value = "BCD0015"

# --- Sheet1: new "Amount" values in column M for rows 2-4 ---
$ws1.Range("M2").Value = 1000
$ws1.Range("M3").Value = 1000
$ws1.Range("M4").Value = 1000

# --- Sheet1: widen the new column R ---
$ws1.Cells.Item(1, 18).ColumnWidth = 27.6

# --- Sheet1: scroll / selection state ---
$ws1.Range("M7").Select()
